# Apply "Added all symbiont density sample data" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 (fragment A16) ----
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 161
$ws.Range("D3").Value = 178
$ws.Range("E3").Value = 182
$ws.Range("F3").Value = 172
$ws.Range("G3").Value = 158
$ws.Range("H3").Value = 185
$ws.Range("I3").Value = "AH, TA"
$ws.Range("N3").Value = 20221017

# ---- Row 4 (fragment A20) ----
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 122
$ws.Range("D4").Value = 112
$ws.Range("E4").Value = 120
$ws.Range("F4").Value = 115
$ws.Range("G4").Value = 114
$ws.Range("H4").Value = 130
$ws.Range("I4").Value = "AH, TA "
$ws.Range("N4").Value = 20221017

# ---- Row 8 (fragment B25) ----
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 112
$ws.Range("D8").Value = 118
$ws.Range("E8").Value = 107
$ws.Range("F8").Value = 105
$ws.Range("G8").Value = 105
$ws.Range("H8").Value = 115
$ws.Range("I8").Value = "AH, TA"
$ws.Range("N8").Value = 20221017

# ---- Row 11 (fragment C19) ----
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 124
$ws.Range("D11").Value = 119
$ws.Range("E11").Value = 113
$ws.Range("F11").Value = 104
$ws.Range("G11").Value = 114
$ws.Range("H11").Value = 107
$ws.Range("I11").Value = "AH, TA"
$ws.Range("N11").Value = 20221017

# ---- Row 12 (fragment C26) ----
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 218
$ws.Range("D12").Value = 220
$ws.Range("E12").Value = 219
$ws.Range("F12").Value = 228
$ws.Range("G12").Value = 238
$ws.Range("H12").Value = 204
$ws.Range("I12").Value = "AH"
$ws.Range("N12").Value = 20221017

# ---- Row 17 (fragment D30) already has counts; only add initials ----
$ws.Range("I17").Value = "AH, TA"

# ---- Row 18 (fragment E5) ----
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 152
$ws.Range("D18").Value = 170
$ws.Range("E18").Value = 161
$ws.Range("F18").Value = 160
$ws.Range("G18").Value = 143
$ws.Range("H18").Value = 150
$ws.Range("I18").Value = "AH"
$ws.Range("N18").Value = 20221017

# ---- Row 20 (fragment E20) ----
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 169
$ws.Range("D20").Value = 165
$ws.Range("E20").Value = 155
$ws.Range("F20").Value = 156
$ws.Range("G20").Value = 159
$ws.Range("H20").Value = 158
$ws.Range("I20").Value = "AH"
$ws.Range("N20").Value = 20221017

# ---- Row 21 (fragment E31) ----
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = 181
$ws.Range("D21").Value = 173
$ws.Range("E21").Value = 205
$ws.Range("F21").Value = 204
$ws.Range("G21").Value = 180
$ws.Range("H21").Value = 174
$ws.Range("I21").Value = "AH"
$ws.Range("N21").Value = 20221017

# ---- Row 23 (fragment F17) ----
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 179
$ws.Range("D23").Value = 175
$ws.Range("E23").Value = 186
$ws.Range("F23").Value = 207
$ws.Range("G23").Value = 187
$ws.Range("H23").Value = 171
$ws.Range("I23").Value = "AH"
$ws.Range("N23").Value = 20221017

# ---- Row 24 (fragment F22) ----
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 257
$ws.Range("D24").Value = 232
$ws.Range("E24").Value = 256
$ws.Range("F24").Value = 217
$ws.Range("G24").Value = 250
$ws.Range("H24").Value = 237
$ws.Range("I24").Value = "AH"
$ws.Range("N24").Value = 20221017

# ---- Row 25 ----
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 118
$ws.Range("D25").Value = 106
$ws.Range("E25").Value = 142
$ws.Range("F25").Value = 114
$ws.Range("G25").Value = 131
$ws.Range("H25").Value = 132
$ws.Range("I25").Value = "AH"
$ws.Range("N25").Value = 20221017

# ---- Row 26 no longer exists; remove its (empty-data) formulas ----
$ws.Rows(26).Delete()

# ---- Restore last-selected cell to match authored state ----
$ws.Range("I23").Select()
